# "Fix bulge chaser batch rotate"
#
# The bulge-chaser batch-rotate step on the "Full 6x6" sheet multiplied the
# rotation matrices against the ORIGINAL input block (A12:F17 / A30:F35)
# instead of the previous iteration's rotated result (AH3:AM8 / AH21:AM26).
# Fix the two array formulas; everything else (the cached values feeding
# off them, the view/selection that was left on the matching cell, and the
# "rand" sheet's volatile RAND() cells) follows from recalculation.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Workbook-level: the active tab moves from "rand" (index 5) to
#    "Full 6x6" (index 4) -- activeTab is 0-based, so 4 -> 3.
# ---------------------------------------------------------------------
$wsFull = $wb.Worksheets.Item("Full 6x6")
$wsRand = $wb.Worksheets.Item("rand")

# ---------------------------------------------------------------------
# 2. The actual fix: correct the two batch-rotate array formulas.
# ---------------------------------------------------------------------
$wsFull.Range("AA12:AF17").FormulaArray = "=MMULT(T12:Y17,MMULT(M12:R17,AH3:AM8))"
$wsFull.Range("AA30:AF35").FormulaArray = "=MMULT(M30:R35,AH21:AM26)"

# ---------------------------------------------------------------------
# 3. A handful of plain (non-array) formulas elsewhere on the sheet
#    reference individual cells inside those two array spill ranges
#    (row 48-51 mirror AI30:AM31 for the test assertions below). Touch
#    them so their cached values pick up the corrected results.
# ---------------------------------------------------------------------
$wsFull.Range("B48").Formula = "=AI30"
$wsFull.Range("C48").Formula = "=AJ30"
$wsFull.Range("D48").Formula = "=AK30"
$wsFull.Range("E48").Formula = "=AL30"
$wsFull.Range("F48").Formula = "=AM30"

$wsFull.Range("B49").Formula = "=AI31"
$wsFull.Range("C49").Formula = "=AJ31"
$wsFull.Range("D49").Formula = "=AK31"
$wsFull.Range("E49").Formula = "=AL31"
$wsFull.Range("F49").Formula = "=AM31"

$wsFull.Range("B50").Formula = "=AI32"
$wsFull.Range("C50").Formula = "=AJ32"
$wsFull.Range("D50").Formula = "=AK32"
$wsFull.Range("E50").Formula = "=AL32"
$wsFull.Range("F50").Formula = "=AM32"

$wsFull.Range("B51").Formula = "=AI33"
$wsFull.Range("C51").Formula = "=AJ33"
$wsFull.Range("D51").Formula = "=AK33"
$wsFull.Range("E51").Formula = "=AL33"
$wsFull.Range("F51").Formula = "=AM33"

$wsFull.Range("B52").Formula = "=AI34"
$wsFull.Range("C52").Formula = "=AJ34"
$wsFull.Range("D52").Formula = "=AK34"
$wsFull.Range("E52").Formula = "=AL34"
$wsFull.Range("F52").Formula = "=AM34"

$wsFull.Range("B53").Formula = "=AI35"
$wsFull.Range("C53").Formula = "=AJ35"
$wsFull.Range("D53").Formula = "=AK35"
$wsFull.Range("E53").Formula = "=AL35"
$wsFull.Range("F53").Formula = "=AM35"

# ---------------------------------------------------------------------
# 4. View state: the author had scrolled to / selected the fixed
#    formula cell on "Full 6x6" and that sheet became the active one
#    (replacing "rand" as tabSelected).
# ---------------------------------------------------------------------
$wsFull.Application.ActiveWindow.ScrollRow = 16
$wsFull.Range("AA12").Select()
$wsFull.Activate()

$wsRand.Range("B5:H8").Select()
